$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.193.20"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.956.13"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.73"
$ws.Range("E5").Value = "  +4.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.20"
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.951.04"
$ws.Range("E7").Value = "  -2.36%  "
$ws.Range("E8").Value = "  -6.59%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  -4.48%  "
$ws.Range("E11").Value = "  -5.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.69"
$ws.Range("E12").Value = "  +18.32%  "
$ws.Range("E13").Value = "  -3.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.64"
$ws.Range("E14").Value = "  -4.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.587.22"
$ws.Range("E15").Value = "  -2.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.954.48"
$ws.Range("E16").Value = "  -2.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.41"
$ws.Range("E17").Value = "  -4.10%  "
$ws.Range("E18").Value = "  -3.10%  "
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("E20").Value = "  -3.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.008.16"
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "423.79"
$ws.Range("E22").Value = "  -4.72%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "97.01"
$ws.Range("E24").Value = "  -7.54%  "
$ws.Range("E25").Value = "  +5.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.27"
$ws.Range("E26").Value = "  -3.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.58"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.69"
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.74"
$ws.Range("E29").Value = "  +16.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.82"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.41"
$ws.Range("E31").Value = "  -3.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.81"
$ws.Range("E32").Value = "  +14.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.82"
$ws.Range("E33").Value = "  +18.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "695.64"
$ws.Range("E34").Value = "  +2.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.39"
$ws.Range("E35").Value = "  -2.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.130"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "64.73"
$ws.Range("E37").Value = "  -3.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.435"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.151"
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0820"
$ws.Range("E40").Value = "  -4.82%  "
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.43"
$ws.Range("E41").Value = "  -4.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.21"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0480"
$ws.Range("E45").Value = "  -3.93%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.148"
$ws.Range("E46").Value = "  -6.55%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.79"
$ws.Range("E48").Value = "  +7.70%  "
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("E51").Value = "  +2.89%  "
